$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "98.506.31"
Set-TextValue "E2" "  +3.37%  "

# Row 3
Set-TextValue "D3" "3.640.94"
Set-TextValue "E3" "  +2.55%  "

# Row 4
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.08%  "

# Row 5
Set-TextValue "D5" "245.95"
Set-TextValue "E5" "  +5.00%  "

# Row 6
Set-TextValue "B6" "XRP"
Set-TextValue "C6" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D6" "1.75"
Set-TextValue "E6" "  +20.77%  "

# Row 7
Set-TextValue "B7" "BNB"
Set-TextValue "C7" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue "D7" "660.75"
Set-TextValue "E7" "  +1.93%  "

# Row 8
Set-TextValue "D8" "0.423"
Set-TextValue "E8" "  +6.69%  "

# Row 9
Set-TextValue "D9" "1.10"
Set-TextValue "E9" "  +11.38%  "

# Row 10
Set-TextValue "E10" "  -0.14%  "

# Row 11
Set-TextValue "D11" "3.635.19"
Set-TextValue "E11" "  +2.52%  "

# Row 12
Set-TextValue "D12" "44.35"
Set-TextValue "E12" "  +5.87%  "

# Row 13
Set-TextValue "D13" "0.207"
Set-TextValue "E13" "  +2.73%  "

# Row 14
Set-TextValue "D14" "6.52"
Set-TextValue "E14" "  +0.20%  "

# Row 15
Set-TextValue "D15" "4.316.46"
Set-TextValue "E15" "  +1.83%  "

# Row 16
Set-TextValue "D16" "98.178.74"
Set-TextValue "E16" "  +3.14%  "

# Row 17
Set-TextValue "D17" "0.0000263"
Set-TextValue "E17" "  +4.75%  "

# Row 18
Set-TextValue "B18" "WrappedEther"
Set-TextValue "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D18" "3.639.50"
Set-TextValue "E18" "  +2.37%  "

# Row 19
Set-TextValue "B19" "Polkadot"
Set-TextValue "C19" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D19" "8.04"
Set-TextValue "E19" "  +2.01%  "

# Row 20
Set-TextValue "D20" "12.94"
Set-TextValue "E20" "  +2.18%  "

# Row 21
Set-TextValue "D21" "18.35"
Set-TextValue "E21" "  +3.91%  "

# Row 22
Set-TextValue "D22" "0.544"
Set-TextValue "E22" "  +15.53%  "

# Row 23
Set-TextValue "B23" "BitcoinCash"
Set-TextValue "C23" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D23" "520.82"
Set-TextValue "E23" "  +3.54%  "

# Row 24
Set-TextValue "B24" "SuiNetwork"
Set-TextValue "C24" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D24" "3.48"
Set-TextValue "E24" "  +0.16%  "

# Row 25
Set-TextValue "E25" "  +9.43%  "

# Row 26
Set-TextValue "D26" "6.96"
Set-TextValue "E26" "  +6.49%  "

# Row 27
Set-TextValue "D27" "99.77"
Set-TextValue "E27" "  +5.25%  "

# Row 28
Set-TextValue "D28" "13.11"
Set-TextValue "E28" "  +5.61%  "

# Row 29
Set-TextValue "D29" "3.834.97"
Set-TextValue "E29" "  +2.39%  "

# Row 30
Set-TextValue "D30" "0.157"
Set-TextValue "E30" "  +14.11%  "

# Row 31
Set-TextValue "D31" "3.07"
Set-TextValue "E31" "  +1.32%  "

# Row 32
Set-TextValue "D32" "11.94"
Set-TextValue "E32" "  +6.90%  "

# Row 33
Set-TextValue "E33" "  -0.02%  "

# Row 34
Set-TextValue "D34" "0.187"
Set-TextValue "E34" "  +5.75%  "

# Row 35
Set-TextValue "D35" "0.995"
Set-TextValue "E35" "  -0.26%  "

# Row 36
Set-TextValue "D36" "32.30"
Set-TextValue "E36" "  +1.71%  "

# Row 37
Set-TextValue "D37" "8.90"
Set-TextValue "E37" "  +9.53%  "

# Row 38
Set-TextValue "B38" "PolygonEcosystemToken"
Set-TextValue "C38" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D38" "0.577"
Set-TextValue "E38" "  +4.03%  "

# Row 39
Set-TextValue "B39" "Bittensor"
Set-TextValue "C39" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D39" "616.45"
Set-TextValue "E39" "  +9.96%  "

# Row 40
Set-TextValue "D40" "1.65"
Set-TextValue "E40" "  +10.10%  "

# Row 41
Set-TextValue "D41" "2.00"
Set-TextValue "E41" "  +13.95%  "

# Row 42
Set-TextValue "E42" "  +3.42%  "

# Row 43
Set-TextValue "D43" "0.938"
Set-TextValue "E43" "  +4.58%  "

# Row 45
Set-TextValue "D45" "6.11"
Set-TextValue "E45" "  +8.98%  "

# Row 46
Set-TextValue "D46" "0.0448"
Set-TextValue "E46" "  +9.42%  "

# Row 47
Set-TextValue "E47" "  +1.18%  "

# Row 48
Set-TextValue "E48" "  +0.66%  "

# Row 49
Set-TextValue "D49" "8.69"
Set-TextValue "E49" "  +8.30%  "

# Row 50
Set-TextValue "D50" "33.13"
Set-TextValue "E50" "  -4.59%  "

# Row 51
Set-TextValue "B51" "Algorand"
Set-TextValue "C51" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D51" "0.397"
Set-TextValue "E51" "  +35.29%  "
